$d = $word.ActiveDocument

# Replace whole-paragraph text while preserving the paragraph's other runs
# (e.g. a leading empty <w:r/> run, or bold/italic formatting on the text
# run) as closely as possible. Word's Find/Replace (and any operation that
# deletes text touching a zero-length run) normalizes/merges the adjacent
# empty run away, so instead we insert the replacement text first (which
# does not disturb the empty run) and only then delete the old text.
function Replace-ParagraphText($oldText, $newText) {
    $matched = 0
    foreach ($p in $d.Paragraphs) {
        $rng = $p.Range
        if ($rng.Text -eq ($oldText + "`r")) {
            # Only re-apply Bold/Italic explicitly if the run carries DIRECT
            # character formatting that differs from what its paragraph
            # style already provides (otherwise we'd bake an explicit
            # <w:b/>/<w:i/> into e.g. a Heading1 run that was only ever
            # bold because of its style).
            $styleName = $rng.ParagraphStyle.NameLocal
            $style = $d.Styles($styleName)
            $directBold = $null
            $directItalic = $null
            if ($rng.Bold -ne $style.Font.Bold) { $directBold = $rng.Bold }
            if ($rng.Italic -ne $style.Font.Italic) { $directItalic = $rng.Italic }

            $insPoint = $d.Range($rng.Start, $rng.Start)
            $insPoint.InsertBefore($newText)

            $newRng = $d.Range($rng.Start, $rng.Start + $newText.Length)
            if ($directBold -ne $null) { $newRng.Bold = $directBold }
            if ($directItalic -ne $null) { $newRng.Italic = $directItalic }

            $oldStart = $rng.Start + $newText.Length
            $oldEnd = $oldStart + $oldText.Length
            $oldRng = $d.Range($oldStart, $oldEnd)
            $oldRng.Delete()

            $matched = $matched + 1
        }
    }
    return $matched
}

# Title heading and the bold "Play Burning Reels..." line near the end (both instances)
Replace-ParagraphText "Play Burning Reels for Free - Unique Firefighter Themed Slot Game" "Play Burning Reels Free - Exciting Gameplay with Firefighter Theme"

# "What we like" bullet list items
Replace-ParagraphText "Special features include scatter with free spins and multiplier symbol" "Special features like Wilds and Scatters with free spins"
Replace-ParagraphText "Customizable volatility levels, animation, screen size, and spin speed with Wazdan tools" "Highly customizable gameplay with Wazdan tools"
Replace-ParagraphText "Highly entertaining with realistic fire animations" "Exciting winning potential with up to 5,000 times line bet"

# "What we don't like" bullet list item
Replace-ParagraphText "No progressive jackpot" "Lack of additional bonus features"

# Italic summary paragraph near the end
Replace-ParagraphText "Join a team of firefighters and prevent a forest fire in Burning Reels. Spin with Wazdan tools for customizable gameplay and enjoy free spins and multipliers. Play now for free." "Play Burning Reels free and experience an action-packed online slot game with a unique firefighter theme."
